$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Update existing rows 22-28 (values shift / change as per new weekly data)
# Columns A,B,C,E,F,G,H,I,J stay the same for all these rows; only
# D,K,L,M,N,O,P,Q,R,S,T change.
# ---------------------------------------------------------------------------

function Set-Row {
    param($r, $d, $k, $l, $m, $n, $o, $p, $q, $rOrigin, $s, $t)

    $ws.Range("D$r").Value = $d
    $ws.Range("K$r").Value = $k
    $ws.Range("L$r").Value = $l
    $ws.Range("M$r").Value = $m
    $ws.Range("N$r").Value = $n
    $ws.Range("O$r").Value = $o
    $ws.Range("P$r").Value = $p
    $ws.Range("Q$r").Value = $q
    $ws.Range("R$r").Value = $rOrigin
    $ws.Range("S$r").Value = $s
    $ws.Range("T$r").Value = $t
}

Set-Row 22 44694 "Wonderfull" "Especial" 200 21600 21600 21600 "`$/caja 18 kilos granel" "Región de O'Higgins" 1200 18
Set-Row 23 44694 "Wonderfull" "Primera"  220 18000 18000 18000 "`$/caja 18 kilos granel" "Región de O'Higgins" 1000 18
Set-Row 24 44694 "Wonderfull" "Segunda"  250 14400 14400 14400 "`$/caja 18 kilos granel" "Región de O'Higgins" 800  18

Set-Row 25 44644 "Sin especificar" "Especial" 180 18000 18000 18000 "`$/caja 15 kilos granel" "Provincia de Limarí" 1200 15
Set-Row 26 44644 "Sin especificar" "Primera"  220 13500 13500 13500 "`$/caja 15 kilos granel" "Provincia de Limarí" 900  15
Set-Row 27 44644 "Sin especificar" "Segunda"  290 12000 12000 12000 "`$/caja 15 kilos granel" "Provincia de Limarí" 800  15

Set-Row 28 44678 "Sin especificar" "Especial" 290 15000 15000 15000 "`$/caja 15 kilos granel" "Región de O'Higgins" 1000 15

# ---------------------------------------------------------------------------
# Append new rows 29-31 (brand new data points, full row content required)
# ---------------------------------------------------------------------------

function Add-Row {
    param($r, $d, $k, $l, $m, $n, $o, $p, $q, $rOrigin, $s, $t)

    $ws.Range("A$r").Value = 9
    $ws.Range("B$r").Value = "Vega Central Mapocho de Santiago"
    $ws.Range("C$r").Value = "Metropolitana"
    $ws.Range("D$r").Value = $d
    $ws.Range("D$r").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range("E$r").Value = 13
    $ws.Range("F$r").Value = "Fruta"
    $ws.Range("G$r").Value = 100104
    $ws.Range("H$r").Value = "Frutos de pepita"
    $ws.Range("I$r").Value = 100104001
    $ws.Range("J$r").Value = "Granada"
    $ws.Range("K$r").Value = $k
    $ws.Range("L$r").Value = $l
    $ws.Range("M$r").Value = $m
    $ws.Range("N$r").Value = $n
    $ws.Range("O$r").Value = $o
    $ws.Range("P$r").Value = $p
    $ws.Range("Q$r").Value = $q
    $ws.Range("R$r").Value = $rOrigin
    $ws.Range("S$r").Value = $s
    $ws.Range("T$r").Value = $t
}

Add-Row 29 44678 "Sin especificar" "Primera" 220 12000 12000 12000 "`$/caja 15 kilos granel" "Región de O'Higgins" 800  15
Add-Row 30 44305 "Wonderfull"      "Primera" 50  18000 18000 18000 "`$/caja 15 kilos granel" "Región de O'Higgins" 1200 15
Add-Row 31 44305 "Wonderfull"      "Segunda" 60  15000 15000 15000 "`$/caja 15 kilos granel" "Región de O'Higgins" 1000 15
